# Swap the order of "System" and the recorded-by email address in column G
# ("Recorded By") wherever the value is exactly "System, dnasr281@gmail.com",
# turning it into "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G ("Recorded By")
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
